$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 28 (shifts the old row 28 down to row 31,
# and leaves 3 empty rows at 28-30 for the new data).
$ws.Range("A28:T30").Insert()

# --- Row 26: update existing "Especial" entry to the new week's data ---
$ws.Range("D26").Value = 44449
$ws.Range("N26").Value = 2900
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = 2950
$ws.Range("Q26").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S26").Value = 2950
$ws.Range("T26").Value = 1

# --- Row 27: update existing entry, quality changes to "Extra (doble especial)" ---
$ws.Range("D27").Value = 44449
$ws.Range("L27").Value = "Extra (doble especial)"
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 3100
$ws.Range("O27").Value = 3200
$ws.Range("P27").Value = 3150
$ws.Range("Q27").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S27").Value = 3150
$ws.Range("T27").Value = 1

# --- Row 28: brand-new "Primera" entry for the new week ---
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44449
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100107
$ws.Range("H28").Value = "Otros"
$ws.Range("I28").Value = 100107002
$ws.Range("J28").Value = "Chirimoya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = 2700
$ws.Range("O28").Value = 2800
$ws.Range("P28").Value = 2750
$ws.Range("Q28").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R28").Value = "Provincia de Limarí"
$ws.Range("S28").Value = 2750
$ws.Range("T28").Value = 1

# --- Rows 29 & 30: new blank rows created by Insert() need to be filled in
#     with what was previously in rows 26 & 27 (the older week's "Especial"
#     and "Primera" entries), now pushed down beneath the new data. ---
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44161
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107002
$ws.Range("J29").Value = "Chirimoya"
$ws.Range("K29").Value = "Cultivar IV Región"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 240
$ws.Range("N29").Value = 13500
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 13750
$ws.Range("Q29").Value = "$/bandeja 8 kilos"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 1719
$ws.Range("T29").Value = 8

$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44161
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107002
$ws.Range("J30").Value = "Chirimoya"
$ws.Range("K30").Value = "Cultivar IV Región"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 11500
$ws.Range("O30").Value = 12000
$ws.Range("P30").Value = 11750
$ws.Range("Q30").Value = "$/bandeja 8 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 1469
$ws.Range("T30").Value = 8

# Row 31 already holds the old row-28 content (unchanged) thanks to Insert(),
# including its date-formatted D cell style.
